$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Info_DE.xlsx" help sheet: column D reports how many characters DeepL
# bills when translating each sample file. Adding the language-code info
# to the help text grew the translated docx/pptx sources slightly, which
# bumped their billed-character counts.
$ws.Range("D3").Value = 1035
$ws.Range("D4").Value = 1034

# Reflect the author's final view state: scrolled one column to the right
# (column B becomes the left-most visible column) with D5 selected.
$ws.Range("D5").Select()
$excel.ActiveWindow.ScrollColumn = 2
